$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 169, shifting existing rows 169:205 down to 170:206
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row 169 with the new data record
$ws.Range("A169").Value = 6
$ws.Range("B169").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 44642
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100112022
$ws.Range("G169").Value = "Arveja Verde"
$ws.Range("H169").Value = "Perfection"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 470
$ws.Range("K169").Value = 25000
$ws.Range("L169").Value = 25000
$ws.Range("M169").Value = 25000
$ws.Range("N169").Value = "`$/saco 25 kilos"
$ws.Range("O169").Value = "Carahue"
$ws.Range("P169").Value = 1000
$ws.Range("Q169").Value = 25
$ws.Range("R169").Value = "Hortaliza"
